$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert first new record row at row 108, pushing existing rows 108.. down by one.
$ws.Rows.Item(108).Insert()

$ws.Cells.Item(108,1).Value2  = 2
$ws.Cells.Item(108,2).Value   = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(108,3).Value   = "Coquimbo"
$ws.Cells.Item(108,4).Value2  = 45210
$ws.Cells.Item(108,5).Value2  = 4
$ws.Cells.Item(108,6).Value2  = 100112024
$ws.Cells.Item(108,7).Value   = "Choclo"
$ws.Cells.Item(108,8).Value   = "Dulce o Americano"
$ws.Cells.Item(108,9).Value   = "Primera"
$ws.Cells.Item(108,10).Value2 = 1000
$ws.Cells.Item(108,11).Value2 = 38000
$ws.Cells.Item(108,12).Value2 = 40000
$ws.Cells.Item(108,13).Value2 = 39000
$ws.Cells.Item(108,14).Value  = "`$/malla 70 unidades"
$ws.Cells.Item(108,15).Value  = "Provincia de Limarí"
$ws.Cells.Item(108,16).Value2 = 557
$ws.Cells.Item(108,17).Value2 = 70
$ws.Cells.Item(108,18).Value  = "Hortaliza"

# Insert second new record row at row 134 (original row 133 position, now shifted to 134
# by the insertion above), pushing the remaining rows down by one more.
$ws.Rows.Item(134).Insert()

$ws.Cells.Item(134,1).Value2  = 2
$ws.Cells.Item(134,2).Value   = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(134,3).Value   = "Coquimbo"
$ws.Cells.Item(134,4).Value2  = 45169
$ws.Cells.Item(134,5).Value2  = 4
$ws.Cells.Item(134,6).Value2  = 100112024
$ws.Cells.Item(134,7).Value   = "Choclo"
$ws.Cells.Item(134,8).Value   = "Dulce o Americano"
$ws.Cells.Item(134,9).Value   = "Primera"
$ws.Cells.Item(134,10).Value2 = 300
$ws.Cells.Item(134,11).Value2 = 45000
$ws.Cells.Item(134,12).Value2 = 46000
$ws.Cells.Item(134,13).Value2 = 45500
$ws.Cells.Item(134,14).Value  = "`$/malla 70 unidades"
$ws.Cells.Item(134,15).Value  = "Provincia de Limarí"
$ws.Cells.Item(134,16).Value2 = 650
$ws.Cells.Item(134,17).Value2 = 70
$ws.Cells.Item(134,18).Value  = "Hortaliza"
